$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new question rows
$ws.Range("A9").Value = "What is your favorite music genre?"
$ws.Range("B9").Value = "String"

$ws.Range("A10").Value = "What is your favorite part of the day?"
$ws.Range("B10").Value = "String"

# Widen column A to fit the longer question text
# (target stored width is 52.26953125; the host quantizes ColumnWidth to
# whole pixels at 6px/char, so 51.5 is the closest input that reproduces it)
$ws.Columns.Item(1).ColumnWidth = 51.5

# Update the selected cell to A11 (next empty row), matching the saved view state
$ws.Range("A11").Select()
